$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row for a5cddbcd file -> zh-cn / de-de status columns
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn / de-de detail sheets: Status column (C) for the a5cddbcd row
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail column (P) for the a5cddbcd row: report the handback/handoff file mismatch
$zhcn.Range("P3").Value = "Handback file name: 4kvvv2tq.vtr is different with handoff file name: a5cddbcd-fd0e-4a54-b20e-0fcea14b3d1e.ae6f3d0645cac3af796e126b59ebad95a695be45.zh-cn."
$dede.Range("P3").Value = "Handback file name: 4kvvv2tq.vtr is different with handoff file name: a5cddbcd-fd0e-4a54-b20e-0fcea14b3d1e.ae6f3d0645cac3af796e126b59ebad95a695be45.de-de."

# Widen the Error Detail column (16th / P) now that it holds a long message
# (ColumnWidth is stored with a fixed +5/6-character padding offset, so request
# 40 - 0.8333333333333334 to land exactly on a stored width of 40)
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
